$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.823.18'
$ws.Range("E2").Value = '  +1.06%  '

$ws.Range("D3").Value = '2.088.64'
$ws.Range("E3").Value = '  +0.98%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '235.21'
$ws.Range("E5").Value = '  +0.20%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '58.82'
$ws.Range("E7").Value = '  +2.91%  '

$ws.Range("E8").Value = '  -0.06%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.393'
$ws.Range("E9").Value = '  +0.22%  '

$ws.Range("E10").Value = '  +1.89%  '

$ws.Range("E11").Value = '  +2.91%  '

$ws.Range("D12").Value = '2.397.14'
$ws.Range("E12").Value = '  +1.10%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.78'
$ws.Range("E13").Value = '  +2.97%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '21.27'
$ws.Range("E14").Value = '  +2.95%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.772'
$ws.Range("E15").Value = '  -0.73%  '

$ws.Range("E16").Value = '  +2.57%  '

$ws.Range("D17").Value = '2.092.58'
$ws.Range("E17").Value = '  +1.19%  '

$ws.Range("D18").Value = '37.753.71'
$ws.Range("E18").Value = '  +1.11%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.21'
$ws.Range("E19").Value = '  +0.03%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '71.33'
$ws.Range("E20").Value = '  +2.61%  '

$ws.Range("E21").Value = '  +2.33%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '228.60'
$ws.Range("E22").Value = '  +1.02%  '

$ws.Range("E23").Value = '  -0.09%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.45'
$ws.Range("E24").Value = '  +2.27%  '

$ws.Range("E25").Value = '  -1.05%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '170.24'
$ws.Range("E26").Value = '  +1.24%  '

$ws.Range("E27").Value = '  +5.07%  '

$ws.Range("E28").Value = '  +2.20%  '

$ws.Range("E29").Value = '  -0.01%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.54'
$ws.Range("E30").Value = '  +2.42%  '

$ws.Range("E31").Value = '  +2.07%  '

$ws.Range("E32").Value = '  +3.04%  '

$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0633'
$ws.Range("E33").Value = '  +2.62%  '

$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.70'
$ws.Range("E34").Value = '  +3.49%  '

$ws.Range("E35").Value = '  +1.57%  '

$ws.Range("E36").Value = '  +3.10%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.83'
$ws.Range("E37").Value = '  +3.31%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.00'
$ws.Range("E38").Value = '  -0.05%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.38'
$ws.Range("E39").Value = '  -4.16%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0996'
$ws.Range("E40").Value = '  +4.66%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '98.86'
$ws.Range("E41").Value = '  +1.93%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.93'
$ws.Range("E42").Value = '  -0.35%  '

$ws.Range("E43").Value = '  +4.41%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0215'
$ws.Range("E44").Value = '  +1.40%  '

$ws.Range("D45").Value = '1.465.40'
$ws.Range("E45").Value = '  -1.80%  '

$ws.Range("E46").Value = '  +0.76%  '

$ws.Range("E47").Value = '  +4.40%  '

$ws.Range("E48").Value = '  +5.29%  '

$ws.Range("E49").Value = '  +2.67%  '

$ws.Range("E50").Value = '  +2.54%  '

$ws.Range("D51").Value = '2.280.91'
